$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 35.88321066666667
$ws.Cells.Item(2, 8).Value = 107.649632
$ws.Cells.Item(2, 9).Value = 0.08317795499144418
$ws.Cells.Item(2, 10).Value = 0.08448843719082051
$ws.Cells.Item(2, 13).Value = 8.850733666666665
$ws.Cells.Item(2, 14).Value = 26.552201
$ws.Cells.Item(2, 15).Value = 0.03425243510433672
$ws.Cells.Item(2, 16).Value = 0.03471479367970558
$ws.Cells.Item(2, 17).Value = 317.5927407155591
$ws.Cells.Item(2, 18).Value = 2858.334666440032
$ws.Cells.Item(2, 19).Value = 0.002849047505455882
$ws.Cells.Item(2, 20).Value = 0.002932998665400097
$ws.Cells.Item(3, 7).Value = 35.88321066666667
$ws.Cells.Item(3, 8).Value = 107.649632
$ws.Cells.Item(3, 9).Value = 0.08317795499144418
$ws.Cells.Item(3, 10).Value = 0.08448843719082051
$ws.Cells.Item(3, 15).Value = 0.06418224226698629
$ws.Cells.Item(3, 16).Value = 0.06504861016194469
$ws.Cells.Item(3, 17).Value = 595.1055498609315
$ws.Cells.Item(3, 18).Value = 5355.949948748384
$ws.Cells.Item(3, 19).Value = 0.005338547658533352
$ws.Cells.Item(3, 20).Value = 0.005495855414017633
$ws.Cells.Item(4, 7).Value = 35.88321066666667
$ws.Cells.Item(4, 8).Value = 107.649632
$ws.Cells.Item(4, 9).Value = 0.08317795499144418
$ws.Cells.Item(4, 10).Value = 0.08448843719082051
$ws.Cells.Item(4, 13).Value = 94.82067633333332
$ws.Cells.Item(4, 14).Value = 284.462029
$ws.Cells.Item(4, 15).Value = 0.3669570438989389
$ws.Cells.Item(4, 16).Value = 0.3719104358409092
$ws.Cells.Item(4, 17).Value = 3402.470304424814
$ws.Cells.Item(4, 18).Value = 30622.23273982333
$ws.Cells.Item(4, 19).Value = 0.03052273648121934
$ws.Cells.Item(4, 20).Value = 0.03142213149915533
$ws.Cells.Item(5, 7).Value = 35.88321066666667
$ws.Cells.Item(5, 8).Value = 107.649632
$ws.Cells.Item(5, 9).Value = 0.08317795499144418
$ws.Cells.Item(5, 10).Value = 0.08448843719082051
$ws.Cells.Item(5, 13).Value = 10.324604
$ws.Cells.Item(5, 14).Value = 20.649208
$ws.Cells.Item(5, 15).Value = 0.03995632924983982
$ws.Cells.Item(5, 16).Value = 0.02699712145781534
$ws.Cells.Item(5, 17).Value = 370.4799403819094
$ws.Cells.Item(5, 18).Value = 2222.879642291457
$ws.Cells.Item(5, 19).Value = 0.003323485755966501
$ws.Cells.Item(5, 20).Value = 0.002280944600621584
$ws.Cells.Item(6, 7).Value = 35.88321066666667
$ws.Cells.Item(6, 8).Value = 107.649632
$ws.Cells.Item(6, 9).Value = 0.08317795499144418
$ws.Cells.Item(6, 10).Value = 0.08448843719082051
$ws.Cells.Item(6, 13).Value = 127.8166836666667
$ws.Cells.Item(6, 14).Value = 383.450051
$ws.Cells.Item(6, 15).Value = 0.4946519494798983
$ws.Cells.Item(6, 16).Value = 0.5013290388596253
$ws.Cells.Item(6, 17).Value = 4586.472986725693
$ws.Cells.Item(6, 18).Value = 41278.25688053123
$ws.Cells.Item(6, 19).Value = 0.0411441375902691
$ws.Cells.Item(6, 20).Value = 0.04235650701162587
$ws.Cells.Item(7, 9).Value = 0.03522729558434242
$ws.Cells.Item(7, 10).Value = 0.03578230735158529
$ws.Cells.Item(7, 13).Value = 8.850733666666665
$ws.Cells.Item(7, 14).Value = 26.552201
$ws.Cells.Item(7, 15).Value = 0.03425243510433672
$ws.Cells.Item(7, 16).Value = 0.03471479367970558
$ws.Cells.Item(7, 17).Value = 134.5059920477635
$ws.Cells.Item(7, 18).Value = 1210.553928429872
$ws.Cells.Item(7, 19).Value = 0.001206620655903976
$ws.Cells.Item(7, 20).Value = 0.001242175417094096
$ws.Cells.Item(8, 9).Value = 0.03522729558434242
$ws.Cells.Item(8, 10).Value = 0.03578230735158529
$ws.Cells.Item(8, 15).Value = 0.06418224226698629
$ws.Cells.Item(8, 16).Value = 0.06504861016194469
$ws.Cells.Item(8, 19).Value = 0.002260966819605002
$ws.Cells.Item(8, 20).Value = 0.002327589361608159
$ws.Cells.Item(9, 9).Value = 0.03522729558434242
$ws.Cells.Item(9, 10).Value = 0.03578230735158529
$ws.Cells.Item(9, 13).Value = 94.82067633333332
$ws.Cells.Item(9, 14).Value = 284.462029
$ws.Cells.Item(9, 15).Value = 0.3669570438989389
$ws.Cells.Item(9, 16).Value = 0.3719104358409092
$ws.Cells.Item(9, 17).Value = 1441.004736690743
$ws.Cells.Item(9, 18).Value = 12969.04263021669
$ws.Cells.Item(9, 19).Value = 0.01292690425218444
$ws.Cells.Item(9, 20).Value = 0.01330781352252145
$ws.Cells.Item(10, 9).Value = 0.03522729558434242
$ws.Cells.Item(10, 10).Value = 0.03578230735158529
$ws.Cells.Item(10, 13).Value = 10.324604
$ws.Cells.Item(10, 14).Value = 20.649208
$ws.Cells.Item(10, 15).Value = 0.03995632924983982
$ws.Cells.Item(10, 16).Value = 0.02699712145781534
$ws.Cells.Item(10, 17).Value = 156.9046313923627
$ws.Cells.Item(10, 18).Value = 941.427788354176
$ws.Cells.Item(10, 19).Value = 0.001407553420949414
$ws.Cells.Item(10, 20).Value = 0.000966019297611627
$ws.Cells.Item(11, 9).Value = 0.03522729558434242
$ws.Cells.Item(11, 10).Value = 0.03578230735158529
$ws.Cells.Item(11, 13).Value = 127.8166836666667
$ws.Cells.Item(11, 14).Value = 383.450051
$ws.Cells.Item(11, 15).Value = 0.4946519494798983
$ws.Cells.Item(11, 16).Value = 0.5013290388596253
$ws.Cells.Item(11, 17).Value = 1942.45025150723
$ws.Cells.Item(11, 18).Value = 17482.05226356507
$ws.Cells.Item(11, 19).Value = 0.01742525043569959
$ws.Cells.Item(11, 20).Value = 0.01793870975274996
$ws.Cells.Item(12, 7).Value = 177.70077
$ws.Cells.Item(12, 8).Value = 533.10231
$ws.Cells.Item(12, 9).Value = 0.4119137160358794
$ws.Cells.Item(12, 10).Value = 0.4184034835782469
$ws.Cells.Item(12, 13).Value = 8.850733666666665
$ws.Cells.Item(12, 14).Value = 26.552201
$ws.Cells.Item(12, 15).Value = 0.03425243510433672
$ws.Cells.Item(12, 16).Value = 0.03471479367970558
$ws.Cells.Item(12, 17).Value = 1572.78218763159
$ws.Cells.Item(12, 18).Value = 14155.03968868431
$ws.Cells.Item(12, 19).Value = 0.01410904782710514
$ws.Cells.Item(12, 20).Value = 0.01452479060728892
$ws.Cells.Item(13, 7).Value = 177.70077
$ws.Cells.Item(13, 8).Value = 533.10231
$ws.Cells.Item(13, 9).Value = 0.4119137160358794
$ws.Cells.Item(13, 10).Value = 0.4184034835782469
$ws.Cells.Item(13, 15).Value = 0.06418224226698629
$ws.Cells.Item(13, 16).Value = 0.06504861016194469
$ws.Cells.Item(13, 17).Value = 2947.08061170783
$ws.Cells.Item(13, 18).Value = 26523.72550537047
$ws.Cells.Item(13, 19).Value = 0.02643754591570941
$ws.Cells.Item(13, 20).Value = 0.02721656509368101
$ws.Cells.Item(14, 7).Value = 177.70077
$ws.Cells.Item(14, 8).Value = 533.10231
$ws.Cells.Item(14, 9).Value = 0.4119137160358794
$ws.Cells.Item(14, 10).Value = 0.4184034835782469
$ws.Cells.Item(14, 13).Value = 94.82067633333332
$ws.Cells.Item(14, 14).Value = 284.462029
$ws.Cells.Item(14, 15).Value = 0.3669570438989389
$ws.Cells.Item(14, 16).Value = 0.3719104358409092
$ws.Cells.Item(14, 17).Value = 16849.70719635411
$ws.Cells.Item(14, 18).Value = 151647.364767187
$ws.Cells.Item(14, 19).Value = 0.1511546395779532
$ws.Cells.Item(14, 20).Value = 0.1556086219349405
$ws.Cells.Item(15, 7).Value = 177.70077
$ws.Cells.Item(15, 8).Value = 533.10231
$ws.Cells.Item(15, 9).Value = 0.4119137160358794
$ws.Cells.Item(15, 10).Value = 0.4184034835782469
$ws.Cells.Item(15, 13).Value = 10.324604
$ws.Cells.Item(15, 14).Value = 20.649208
$ws.Cells.Item(15, 15).Value = 0.03995632924983982
$ws.Cells.Item(15, 16).Value = 0.02699712145781534
$ws.Cells.Item(15, 17).Value = 1834.69008074508
$ws.Cells.Item(15, 18).Value = 11008.14048447048
$ws.Cells.Item(15, 19).Value = 0.01645856006045462
$ws.Cells.Item(15, 20).Value = 0.01129568966453498
$ws.Cells.Item(16, 7).Value = 177.70077
$ws.Cells.Item(16, 8).Value = 533.10231
$ws.Cells.Item(16, 9).Value = 0.4119137160358794
$ws.Cells.Item(16, 10).Value = 0.4184034835782469
$ws.Cells.Item(16, 13).Value = 127.8166836666667
$ws.Cells.Item(16, 14).Value = 383.450051
$ws.Cells.Item(16, 15).Value = 0.4946519494798983
$ws.Cells.Item(16, 16).Value = 0.5013290388596253
$ws.Cells.Item(16, 17).Value = 22713.12310641309
$ws.Cells.Item(16, 18).Value = 204418.1079577178
$ws.Cells.Item(16, 19).Value = 0.203753922654657
$ws.Cells.Item(16, 20).Value = 0.2097578162778015
$ws.Cells.Item(17, 7).Value = 20.074196
$ws.Cells.Item(17, 8).Value = 40.148392
$ws.Cells.Item(17, 9).Value = 0.04653236263856699
$ws.Cells.Item(17, 10).Value = 0.0315103250497358
$ws.Cells.Item(17, 13).Value = 8.850733666666665
$ws.Cells.Item(17, 14).Value = 26.552201
$ws.Cells.Item(17, 15).Value = 0.03425243510433672
$ws.Cells.Item(17, 16).Value = 0.03471479367970558
$ws.Cells.Item(17, 17).Value = 177.6713623684653
$ws.Cells.Item(17, 18).Value = 1066.028174210792
$ws.Cells.Item(17, 19).Value = 0.001593846731528978
$ws.Cells.Item(17, 20).Value = 0.001093874432882037
$ws.Cells.Item(18, 7).Value = 20.074196
$ws.Cells.Item(18, 8).Value = 40.148392
$ws.Cells.Item(18, 9).Value = 0.04653236263856699
$ws.Cells.Item(18, 10).Value = 0.0315103250497358
$ws.Cells.Item(18, 15).Value = 0.06418224226698629
$ws.Cells.Item(18, 16).Value = 0.06504861016194469
$ws.Cells.Item(18, 17).Value = 332.9207511437506
$ws.Cells.Item(18, 18).Value = 1997.524506862504
$ws.Cells.Item(18, 19).Value = 0.002986551372123768
$ws.Cells.Item(18, 20).Value = 0.002049702850236424
$ws.Cells.Item(19, 7).Value = 20.074196
$ws.Cells.Item(19, 8).Value = 40.148392
$ws.Cells.Item(19, 9).Value = 0.04653236263856699
$ws.Cells.Item(19, 10).Value = 0.0315103250497358
$ws.Cells.Item(19, 13).Value = 94.82067633333332
$ws.Cells.Item(19, 14).Value = 284.462029
$ws.Cells.Item(19, 15).Value = 0.3669570438989389
$ws.Cells.Item(19, 16).Value = 0.3719104358409092
$ws.Cells.Item(19, 17).Value = 1903.448841567895
$ws.Cells.Item(19, 18).Value = 11420.69304940737
$ws.Cells.Item(19, 19).Value = 0.01707537823948197
$ws.Cells.Item(19, 20).Value = 0.01171901872273596
$ws.Cells.Item(20, 7).Value = 20.074196
$ws.Cells.Item(20, 8).Value = 40.148392
$ws.Cells.Item(20, 9).Value = 0.04653236263856699
$ws.Cells.Item(20, 10).Value = 0.0315103250497358
$ws.Cells.Item(20, 13).Value = 10.324604
$ws.Cells.Item(20, 14).Value = 20.649208
$ws.Cells.Item(20, 15).Value = 0.03995632924983982
$ws.Cells.Item(20, 16).Value = 0.02699712145781534
$ws.Cells.Item(20, 17).Value = 207.258124318384
$ws.Cells.Item(20, 18).Value = 829.0324972735361
$ws.Cells.Item(20, 19).Value = 0.001859262402359528
$ws.Cells.Item(20, 20).Value = 0.0008506880725429586
$ws.Cells.Item(21, 7).Value = 20.074196
$ws.Cells.Item(21, 8).Value = 40.148392
$ws.Cells.Item(21, 9).Value = 0.04653236263856699
$ws.Cells.Item(21, 10).Value = 0.0315103250497358
$ws.Cells.Item(21, 13).Value = 127.8166836666667
$ws.Cells.Item(21, 14).Value = 383.450051
$ws.Cells.Item(21, 15).Value = 0.4946519494798983
$ws.Cells.Item(21, 16).Value = 0.5013290388596253
$ws.Cells.Item(21, 17).Value = 2565.817159994665
$ws.Cells.Item(21, 18).Value = 15394.90295996799
$ws.Cells.Item(21, 19).Value = 0.02301732389307275
$ws.Cells.Item(21, 20).Value = 0.01579704097133842
$ws.Cells.Item(22, 7).Value = 182.547562
$ws.Cells.Item(22, 8).Value = 547.642686
$ws.Cells.Item(22, 9).Value = 0.423148670749767
$ws.Cells.Item(22, 10).Value = 0.4298154468296114
$ws.Cells.Item(22, 13).Value = 8.850733666666665
$ws.Cells.Item(22, 14).Value = 26.552201
$ws.Cells.Item(22, 15).Value = 0.03425243510433672
$ws.Cells.Item(22, 16).Value = 0.03471479367970558
$ws.Cells.Item(22, 17).Value = 1615.67985276132
$ws.Cells.Item(22, 18).Value = 14541.11867485188
$ws.Cells.Item(22, 19).Value = 0.01449387238434274
$ws.Cells.Item(22, 20).Value = 0.01492095455704042
$ws.Cells.Item(23, 7).Value = 182.547562
$ws.Cells.Item(23, 8).Value = 547.642686
$ws.Cells.Item(23, 9).Value = 0.423148670749767
$ws.Cells.Item(23, 10).Value = 0.4298154468296114
$ws.Cells.Item(23, 15).Value = 0.06418224226698629
$ws.Cells.Item(23, 16).Value = 0.06504861016194469
$ws.Cells.Item(23, 17).Value = 3027.462293408931
$ws.Cells.Item(23, 18).Value = 27247.16064068038
$ws.Cells.Item(23, 19).Value = 0.02715863050101476
$ws.Cells.Item(23, 20).Value = 0.02795889744240146
$ws.Cells.Item(24, 7).Value = 182.547562
$ws.Cells.Item(24, 8).Value = 547.642686
$ws.Cells.Item(24, 9).Value = 0.423148670749767
$ws.Cells.Item(24, 10).Value = 0.4298154468296114
$ws.Cells.Item(24, 13).Value = 94.82067633333332
$ws.Cells.Item(24, 14).Value = 284.462029
$ws.Cells.Item(24, 15).Value = 0.3669570438989389
$ws.Cells.Item(24, 16).Value = 0.3719104358409092
$ws.Cells.Item(24, 17).Value = 17309.2832918411
$ws.Cells.Item(24, 18).Value = 155783.5496265699
$ws.Cells.Item(24, 19).Value = 0.1552773853480999
$ws.Cells.Item(24, 20).Value = 0.1598528501615559
$ws.Cells.Item(25, 7).Value = 182.547562
$ws.Cells.Item(25, 8).Value = 547.642686
$ws.Cells.Item(25, 9).Value = 0.423148670749767
$ws.Cells.Item(25, 10).Value = 0.4298154468296114
$ws.Cells.Item(25, 13).Value = 10.324604
$ws.Cells.Item(25, 14).Value = 20.649208
$ws.Cells.Item(25, 15).Value = 0.03995632924983982
$ws.Cells.Item(25, 16).Value = 0.02699712145781534
$ws.Cells.Item(25, 17).Value = 1884.731288815448
$ws.Cells.Item(25, 18).Value = 11308.38773289269
$ws.Cells.Item(25, 19).Value = 0.01690746761010976
$ws.Cells.Item(25, 20).Value = 0.01160377982250419
$ws.Cells.Item(26, 7).Value = 182.547562
$ws.Cells.Item(26, 8).Value = 547.642686
$ws.Cells.Item(26, 9).Value = 0.423148670749767
$ws.Cells.Item(26, 10).Value = 0.4298154468296114
$ws.Cells.Item(26, 13).Value = 127.8166836666667
$ws.Cells.Item(26, 14).Value = 383.450051
$ws.Cells.Item(26, 15).Value = 0.4946519494798983
$ws.Cells.Item(26, 16).Value = 0.5013290388596253
$ws.Cells.Item(26, 17).Value = 23332.62398627522
$ws.Cells.Item(26, 18).Value = 209993.615876477
$ws.Cells.Item(26, 19).Value = 0.2093113149061999
$ws.Cells.Item(26, 20).Value = 0.2154789648461095
